# Passasjertall tom august 2023
# Fill in August figures (column I) on the "Total" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Total")

$ws.Range("I7").Value = 0.67
$ws.Range("I8").Value = 0.57
$ws.Range("I9").Value = 0.56
$ws.Range("I10").Value = 0.61
$ws.Range("I11").Value = 0.7
$ws.Range("I12").Value = 0.66
$ws.Range("I13").Value = 0.73
$ws.Range("I14").Value = 0.5
$ws.Range("I16").Value = 0.63
$ws.Range("I17").Value = 0.8

# Update active selection to match the saved workbook state.
$ws.Range("I18").Select()
